$d = $word.ActiveDocument

# --- 1. Empty paragraph right after the last existing paragraph ---
$anchor = $d.Paragraphs.Last
$anchor.Range.InsertParagraphAfter()
$pEmpty1 = $d.Paragraphs.Last

# --- 2. "Use Cases" heading: centered, bold, dark blue, 12pt ---
$pEmpty1.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Last
$pHeading.Range.Text = "Use Cases"
$pHeading.Range.ParagraphFormat.Alignment = 1
$pHeading.Range.Font.Bold = $true
$pHeading.Range.Font.BoldBi = $true
$pHeading.Range.Font.Color = 6299648
$pHeading.Range.Font.Size = 12
$pHeading.Range.Font.SizeBi = 12

# --- 3-5. Bulleted list of three use cases ---
$pHeading.Range.InsertParagraphAfter()
$pB1 = $d.Paragraphs.Last
$pB1.Range.Text = "Deploy web services or websites"

$pB1.Range.InsertParagraphAfter()
$pB2 = $d.Paragraphs.Last
$pB2.Range.Text = "Lightweight data transformation"

$pB2.Range.InsertParagraphAfter()
$pB3 = $d.Paragraphs.Last
$pB3.Range.Text = "Scheduled tasks or workflows using webhooks"

$bulletsRange = $d.Range($pB1.Range.Start, $pB3.Range.End)
$bulletsRange.Style = $d.Styles.Item("List Paragraph")
$bulletsRange.ListFormat.ApplyBulletDefault()

# --- 6. Empty paragraph carrying the same bold/blue/12pt formatting mark ---
$pB3.Range.InsertParagraphAfter()
$pMark = $d.Paragraphs.Last
$pMark.Range.Style = $d.Styles.Item("Normal")
$pMark.Range.Text = "x"
$pMark.Range.Font.Bold = $true
$pMark.Range.Font.BoldBi = $true
$pMark.Range.Font.Color = 6299648
$pMark.Range.Font.Size = 12
$pMark.Range.Font.SizeBi = 12
$delRange = $d.Range($pMark.Range.Start, $pMark.Range.Start + 1)
$delRange.Delete()

# --- 7-8. Two trailing empty paragraphs ---
$pMark.Range.InsertParagraphAfter()
$pEmpty2 = $d.Paragraphs.Last
$pEmpty2.Range.Style = $d.Styles.Item("Normal")

$pEmpty2.Range.InsertParagraphAfter()
$pEmpty3 = $d.Paragraphs.Last
$pEmpty3.Range.Style = $d.Styles.Item("Normal")

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
